$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.554.02'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '1.878.95'
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.32'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4597'
$ws.Range("E7").Value = '  -1.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3857'
$ws.Range("E8").Value = '  -2.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.65'
$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07854'
$ws.Range("E10").Value = '  -1.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9997'
$ws.Range("E11").Value = '  +2.20%  '

$ws.Range("E12").Value = '  -3.16%  '

$ws.Range("D13").Value = '1.875.02'
$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.051'
$ws.Range("E14").Value = '  +0.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.699'
$ws.Range("E15").Value = '  -0.91%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06972'
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.44'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001002'
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.16'
$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").Value = '28.567.20'
$ws.Range("E22").Value = '  -1.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.329'
$ws.Range("E23").Value = '  -0.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.97'
$ws.Range("E24").Value = '  -1.48%  '

$ws.Range("D25").Value = '2.107.82'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.058'
$ws.Range("E26").Value = '  -2.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.39'
$ws.Range("E27").Value = '  +0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.44'
$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.848'
$ws.Range("E29").Value = '  +1.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.952'
$ws.Range("E30").Value = '  -2.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.33'
$ws.Range("E31").Value = '  -1.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09337'
$ws.Range("E32").Value = '  -0.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9213'
$ws.Range("E33").Value = '  -2.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.309'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.335'
$ws.Range("E35").Value = '  -1.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.272'
$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05777'
$ws.Range("E37").Value = '  -2.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.166'
$ws.Range("E38").Value = '  +1.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.058'
$ws.Range("E39").Value = '  +1.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02067'
$ws.Range("E40").Value = '  -2.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5668'
$ws.Range("E41").Value = '  -1.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1792'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.726'
$ws.Range("E43").Value = '  -2.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.76'
$ws.Range("E44").Value = '  -1.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5355'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07136'
$ws.Range("E46").Value = '  -1.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.150'
$ws.Range("E47").Value = '  +1.10%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.120'
$ws.Range("E48").Value = '  -2.60%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.838'
$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.29'
$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("E51").Value = '  +5.21%  '
